$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.168.38"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.267.23"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D5").Value = "574.88"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("D6").Value = "180.70"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.260.74"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -7.16%  "
$ws.Range("E11").Value = "  -4.41%  "
$ws.Range("D12").Value = "45.64"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("D14").Value = "3.781.57"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "613.53"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "8.29"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.341.14"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "17.48"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.259.44"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "10.74"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "0.876"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").Value = "18.03"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "98.49"
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "4.91"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "30.22"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "8.23"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("D31").Value = "6.40"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "547.58"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  -12.30%  "
$ws.Range("D34").Value = "10.69"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").Value = "3.757.41"
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").Value = "0.102"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "55.49"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  +5.00%  "
$ws.Range("D41").Value = "32.05"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").Value = "3.09"
$ws.Range("E42").Value = "  -6.76%  "
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("D44").Value = "0.0₃0664"
$ws.Range("E44").Value = "  -9.26%  "
$ws.Range("D45").Value = "0.325"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").Value = "2.99"
$ws.Range("E47").Value = "  -7.07%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "0.125"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("D50").Value = "2.47"
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("D51").Value = "127.83"
$ws.Range("E51").Value = "  +4.44%  "
